# Updated Communities overview image
# Updated the Powerpoint and image for the Communities overview
#
# The "Dimension: Communities" slide contains a single table (graphicFrame)
# that lists, for the "Non-Users" category, the Level-2 sub-categories
# Third-Party / Media / Legal-Political, each of which is broken down into
# "Specified" / "Unspecified" columns. Rename those two column headers to
# "Other-Specified" / "Other-Unspecified".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

$specifiedCell = $tbl.Cell(5, 6)
$specifiedCell.Shape.TextFrame.TextRange.Text = "Other-Specified"

$unspecifiedCell = $tbl.Cell(5, 7)
$unspecifiedCell.Shape.TextFrame.TextRange.Text = "Other-Unspecified"
